# feat(modal): modal uses dedicated class - skills about object started - about-me -> about
#
# Moves three overlapping "About"/"Skills" mini-diagram shapes (they were
# stacked on top of each other at the "About" slot) so the "Skills" copy
# lands at the position formerly used by "About", and renames its label
# from "Background" to "Skills".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Reposition "Rectangle 24" (top-level shape #22, id 25) ---------------
$rect24 = $s.Shapes.Item(22)
$rect24.Left = 43.21598434448242
$rect24.Top  = 199.4395294189453

# --- Reposition "Rectangle 27" (top-level shape #25, id 28) ---------------
$rect27 = $s.Shapes.Item(25)
$rect27.Left = 47.35441207885742
$rect27.Top  = 205.4236297607422

# --- Reposition "Group 38" (top-level shape #29, id 39) -------------------
$group38 = $s.Shapes.Item(29)
$group38.Left = 143.5718994140625
$group38.Top  = 199.4395294189453

# --- Rename the "Background" label inside Group 38 to "Skills" ------------
$label = $group38.GroupItems.Item(3)
$label.TextFrame.TextRange.Text = "Skills "
